# "adloori to davuluri completed"
# Grader filled in the "Points for grading" (column E) scores for the
# "Generic" and "Customer Class" rubric sections, matching the max points
# already recorded in column D (i.e. full marks awarded), then left the
# selection on the next section's total cell (E15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Generic section (rows 3-6): award full points in column E ---
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# --- Customer Class section (rows 10-14): award full points in column E ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the selection where the grader ended up after finishing this pass.
$ws.Range("E15").Select() | Out-Null
